# Agregando soporte para apoyos inclinados
# Adds a "rotación" column (D) to the "restric" sheet so that inclined
# supports can also define a rotation/skew angle, with a "grados" comment
# on the header documenting the unit, and switches the active sheet/tab
# from "config" to "restric".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("restric")

# New header cell D1 = "rotación", styled like the existing headers
# (bold, centered, Arial 10) and documented with a "grados" comment,
# matching the existing B1/C1 comments.
$ws.Range("D1").Value = "rotación"
$ws.Range("D1").Font.Name = "Arial"
$ws.Range("D1").Font.Size = 10
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").AddComment("grados")

# New data column values (default rotation = 0) for the existing rows.
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 0

# Make "restric" the active sheet/tab and restore its previous selection,
# moved to the new relevant cell.
$ws.Activate()
$ws.Range("F21").Select()
